# Insert a new data row at row 243 (pushing the existing rows 243-282
# down to 244-283) on the active sheet, then populate the new row with
# the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 243..282 down by one, creating a blank row 243.
$ws.Rows(243).Insert()

# Populate the newly inserted row 243 with the new record.
$ws.Range("A243").Value = 4
$ws.Range("B243").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C243").Value = "Los Lagos"
$ws.Range("D243").Value = 44694
$ws.Range("E243").Value = 10
$ws.Range("F243").Value = 100112040
$ws.Range("G243").Value = "Cilantro"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 200
$ws.Range("K243").Value = 13000
$ws.Range("L243").Value = 13000
$ws.Range("M243").Value = 13000
$ws.Range("N243").Value = "$/caja 36 atados"
$ws.Range("O243").Value = "Región Metropolitana"
$ws.Range("P243").Value = 361
$ws.Range("Q243").Value = 36
$ws.Range("R243").Value = "Hortaliza"
